$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Fix the typo "surpremacy" -> "supremacy" in the "Boss Levels" bullet
# ---------------------------------------------------------------------
$d.Content.Find.Execute("surpremacy", $true, $false, $false, $false, $false, `
    $true, 1, $false, "supremacy", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Insert two new bullets right after the "Adds a new feature..." bullet
#    and before "Breakdown and Tracking:"
#       * Custom Level Creator:                         (ilvl 1, bold label)
#         - Creates a level from a song and a theme      (ilvl 2)
#           for the user to play in. The level is
#           stored locally.
# ---------------------------------------------------------------------
$bossDetail = $d.Paragraphs(22)
$newPara = $bossDetail.Range.InsertParagraphAfter()

# --- "Custom Level Creator:" paragraph (list level 1 / ilvl=1) ---
$creatorPara = $d.Paragraphs(23)
$creatorPara.Range.ListFormat.ListLevelNumber = 2

$r = $creatorPara.Range
$start = $r.Start
$r.InsertAfter("Custom Level Creator")
$boldRange = $d.Range($start, $start + 21)
$boldRange.Font.Bold = 1
$creatorPara.Range.InsertAfter(":")

# --- insert paragraph mark after it, for the description bullet ---
$newPara2 = $creatorPara.Range.InsertParagraphAfter()

# --- description paragraph (list level 2 / ilvl=2) ---
$descPara = $d.Paragraphs(24)
$descPara.Range.ListFormat.ListLevelNumber = 3
$descPara.Range.InsertAfter("Creates a level from a song and a theme")
$descPara.Range.InsertAfter(" for the user to play in.")
$descPara.Range.InsertAfter(" ")
$descPara.Range.InsertAfter("The level is stored locally.")

# ---------------------------------------------------------------------
# 3) "Breakdown and Tracking" -> "Breakdown of Components" and drop the
#    trailing ":" run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Breakdown and Tracking", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Breakdown of Components", 2) | Out-Null

$breakdownPara = $d.Paragraphs(26)
$br = $breakdownPara.Range
$lastChar = $br.Characters($br.Characters.Count - 1)
$lastChar.Delete()

# ---------------------------------------------------------------------
# 4) Move the _GoBack bookmark from the end of "Extra Features: 20 Hours"
#    to the end of "Breakdown of Components".
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$target = $d.Paragraphs(26).Range
$insertPos = $target.End
$target.InsertAfter("Z")
$tempRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $tempRange)
$clearRange = $d.Range($insertPos, $insertPos + 1)
$clearRange.Text = ""

Write-Output "Done"
